$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "7"
${row2} = @{
    "D2" = 0.0413
    "E2" = 0.0287
    "F2" = -0.008920000000000001
    "I2" = 0
    "J2" = 0
    "K2" = 30922.2
    "L2" = 0.2726647467449214
    "M2" = 29186.703
    "N2" = 0.07454159783751314
    "O2" = 0.9438753710926131
    "P2" = 16884.503
    "Q2" = 0.04312230238243368
    "R2" = 0.5460317506516355
    "S2" = 12302.2
    "T2" = 0.4215001605354329
    "U2" = 186859.7
    "V2" = 0.4772317246466088
    "W2" = 0.1086610554624137
    "X2" = 0.0756258827179912
    "Y2" = 0.03303517274442253
    "Z2" = 0.1093498534045842
    "AA2" = 0
    "AB2" = 0.03147083456959962
    "AC2" = -0.03147083456959962
    "AD2" = 1048750.8
    "AE2" = 0
    "AF2" = 1048750.8
    "AG2" = 861891.1000000001
    "AH2" = 0.7281474692772339
    "AI2" = 0.7905695880869865
    "AJ2" = 0.6876203836752337
    "AK2" = 0.7562328877400653
}
foreach ($key in ${row2}.Keys) { $ws.Range($key).Value = ${row2}[$key] }
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# Row 3
$ws.Range("B3").Value = "The Toronto-Dominion Bank (TSX:TD)"
${row3} = @{
    "D3" = 0.0413
    "E3" = 0.08500000000000001
    "F3" = -0.0136
    "I3" = 0
    "J3" = 0
    "K3" = 8939.200000000001
    "L3" = 0.3267502257117688
    "M3" = 11462.438
    "N3" = 0.1118475410971803
    "O3" = 1.282266645784858
    "P3" = 4248.738
    "Q3" = 0.0414580997573249
    "R3" = 0.4752928673706819
    "S3" = 7213.7
    "T3" = 0.6293338293302001
    "U3" = 3420.1
    "V3" = 0.0333724618886895
    "W3" = 0.1436131416178006
    "X3" = 0.0756258827179912
    "Y3" = 0.06798725889980944
    "Z3" = 0.1063118159731775
    "AA3" = 0
    "AB3" = 0.03085305636826789
    "AC3" = -0.03085305636826789
    "AD3" = 269352.8
    "AE3" = 0
    "AF3" = 269352.8
    "AG3" = 265932.7
    "AH3" = 0.7243869937109286
    "AI3" = 0.789610956115559
    "AJ3" = 0.7218284034815048
    "AK3" = 0.7874802191762649
}
foreach ($key in ${row3}.Keys) { $ws.Range($key).Value = ${row3}[$key] }
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# Row 4
$ws.Range("B4").Value = "Royal Bank of Canada (TSX:RY)"
${row4} = @{
    "D4" = 0.046
    "E4" = 0.0287
    "F4" = 0.0737
    "I4" = 0
    "J4" = 0
    "K4" = 8591.200000000001
    "L4" = 0.2669152142169199
    "M4" = 8816.700000000001
    "N4" = 0.07552812376856786
    "O4" = 1.026247788434677
    "P4" = 4557.9
    "Q4" = 0.03904517963917967
    "R4" = 0.5305312412701368
    "S4" = 4258.800000000001
    "T4" = 0.4830378713123959
    "U4" = 87466.39999999999
    "V4" = 0.7492795586547192
    "W4" = 0.1452784339853322
    "X4" = 0.07556993866124739
    "Y4" = 0.06970849532408477
    "Z4" = 0.1080308394869215
    "AA4" = 0
    "AB4" = 0.03085343422783496
    "AC4" = -0.03085343422783496
    "AD4" = 306416.4
    "AE4" = 0
    "AF4" = 306416.4
    "AG4" = 218950
    "AH4" = 0.7241311836169835
    "AI4" = 0.8245369493335171
    "AJ4" = 0.6522503306681284
    "AK4" = 0.7705274567491097
}
foreach ($key in ${row4}.Keys) { $ws.Range($key).Value = ${row4}[$key] }
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# Row 5
$ws.Range("B5").Value = "Canadian Imperial Bank of Commerce (TSX:CM)"
${row5} = @{
    "D5" = 0.0443
    "E5" = 0.0117
    "F5" = 0.0616
    "I5" = 0
    "J5" = 0
    "K5" = 2848.2
    "L5" = 0.2332009661440209
    "M5" = 2128.853
    "N5" = 0.05582617429248747
    "O5" = 0.7474380310371463
    "P5" = 1952.953
    "Q5" = 0.05121344431157824
    "R5" = 0.6856797275472228
    "S5" = 175.9000000000001
    "T5" = 0.08262665388357021
    "U5" = 32365.9
    "V5" = 0.84875018356515
    "W5" = 0.1053694701560454
    "X5" = 0.0869684819006145
    "Y5" = 0.01840098825543086
    "Z5" = 0.09575715856008707
    "AA5" = 0
    "AB5" = 0.03145686969503046
    "AC5" = -0.03145686969503046
    "AD5" = 126238.1
    "AE5" = 0
    "AF5" = 126238.1
    "AG5" = 93872.20000000001
    "AH5" = 0.768003859545165
    "AI5" = 0.802522669826626
    "AJ5" = 0.7111217840428223
    "AK5" = 0.7513641016939113
}
foreach ($key in ${row5}.Keys) { $ws.Range($key).Value = ${row5}[$key] }
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

# Row 6
${row6} = @{
    "D6" = 0.0344
    "E6" = 0.0313
    "F6" = -0.008920000000000001
    "I6" = 0
    "J6" = 0
    "K6" = 3830.4
    "L6" = 0.2292513301772173
    "M6" = 2107.632
    "N6" = 0.0431685691053247
    "O6" = 0.5502380952380951
    "P6" = 2050.532
    "Q6" = 0.04199904553768384
    "R6" = 0.5353310359231411
    "S6" = 57.09999999999991
    "T6" = 0.02709201606352528
    "U6" = 43059.1
    "V6" = 0.881937517537733
    "W6" = 0.1086610554624137
    "X6" = 0.08050218772438414
    "Y6" = 0.02815886773802959
    "Z6" = 0.1240779535376954
    "AA6" = 0
    "AB6" = 0.03147083456959962
    "AC6" = -0.03147083456959962
    "AD6" = 142639
    "AE6" = 0
    "AF6" = 142639
    "AG6" = 99579.89999999999
    "AH6" = 0.7449978403058984
    "AI6" = 0.7703179257867138
    "AJ6" = 0.6710091157063998
    "AK6" = 0.7007245800609246
}
foreach ($key in ${row6}.Keys) { $ws.Range($key).Value = ${row6}[$key] }
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()

# Row 7
$ws.Range("B7").Value = "The Bank of Nova Scotia (TSX:BNS)"
${row7} = @{
    "D7" = 0.027
    "E7" = -0.006820000000000001
    "F7" = -0.095
    "I7" = 0
    "J7" = 0
    "K7" = 5093.7
    "L7" = 0.2684143963745587
    "M7" = 3589.9
    "N7" = 0.0548706300974554
    "O7" = 0.7047725621846596
    "P7" = 3278.8
    "Q7" = 0.05011555251219721
    "R7" = 0.6436971160453109
    "S7" = 311.0999999999999
    "T7" = 0.08665979553748013
    "U7" = 3007.5
    "V7" = 0.04596880693559629
    "W7" = 0.1028052095884522
    "X7" = 0.07199237170660072
    "Y7" = 0.03081283788185148
    "Z7" = 0.1039697749668263
    "AA7" = 0
    "AB7" = 0.03149406715291087
    "AC7" = -0.03149406715291087
    "AD7" = 157657.9
    "AE7" = 0
    "AF7" = 157657.9
    "AG7" = 154650.4
    "AH7" = 0.7067240086299833
    "AI7" = 0.7484658761288142
    "AJ7" = 0.7027161624753719
    "AK7" = 0.7448224976749943
}
foreach ($key in ${row7}.Keys) { $ws.Range($key).Value = ${row7}[$key] }
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()

# Row 8
$ws.Range("B8").Value = "National Bank of Canada (TSX:NA)"
${row8} = @{
    "D8" = 0.0512
    "E8" = 0.0567
    "I8" = 0
    "J8" = 0
    "K8" = 1533.8
    "L8" = 0.2882324200398392
    "M8" = 1001.28
    "N8" = 0.0529999311881686
    "O8" = 0.6528100143434606
    "P8" = 715.6799999999999
    "Q8" = 0.03788250115127487
    "R8" = 0.4666058156213326
    "S8" = 285.6
    "T8" = 0.2852348993288591
    "U8" = 17488.3
    "V8" = 0.9256938085231393
    "W8" = 0.1637170975385863
    "X8" = 0.06404988610718834
    "Y8" = 0.09966721143139798
    "Z8" = 0.2028482773868427
    "AA8" = 0
    "AB8" = 0.03206121084487051
    "AC8" = -0.03206121084487051
    "AD8" = 36501.4
    "AE8" = 0
    "AF8" = 36501.4
    "AG8" = 19013.1
    "AH8" = 0.6589473494182531
    "AI8" = 0.7477757086695634
    "AJ8" = 0.5015960870804007
    "AK8" = 0.6069624900239426
}
foreach ($key in ${row8}.Keys) { $ws.Range($key).Value = ${row8}[$key] }
$ws.Range("F8").ClearContents()
$ws.Range("AN8").ClearContents()
$ws.Range("AP8").ClearContents()

# Row 9
$ws.Range("A9").Value = "Canada"
$ws.Range("B9").Value = "Laurentian Bank of Canada (TSX:LB)"
$ws.Range("C9").Value = "Bank (Money Center)"
${row9} = @{
    "D9" = -0.00175
    "E9" = 0.0217
    "G9" = 0
    "H9" = 0
    "I9" = 0
    "J9" = 0
    "K9" = 85.7
    "L9" = 0.1334267476257201
    "M9" = 79.90000000000001
    "N9" = 0.07546991593463682
    "O9" = 0.9323220536756126
    "P9" = 79.90000000000001
    "Q9" = 0.07546991593463682
    "R9" = 0.9323220536756126
    "S9" = 0
    "T9" = 0
    "U9" = 52.4
    "V9" = 0.04949466326626995
    "W9" = 0.04853324272284517
    "X9" = 0.188122399749355
    "Y9" = -0.1395891570265098
    "Z9" = 0.05912676860196445
    "AA9" = 0
    "AB9" = 0.03563270154781673
    "AC9" = -0.03563270154781673
    "AD9" = 9945.200000000001
    "AE9" = 0
    "AF9" = 9945.200000000001
    "AG9" = 9892.800000000001
    "AH9" = 0.9037886567489708
    "AI9" = 0.8351976888709732
    "AJ9" = 0.9033283111902478
    "AK9" = 0.8344692624333626
    "AL9" = 0
    "AM9" = 0
}
foreach ($key in ${row9}.Keys) { $ws.Range($key).Value = ${row9}[$key] }
